{"js": "// Replace specific three-digit-by-one-digit multiplication answers\n// with their new values, one-to-one, as described by the diff.\nconst replacements = [\n  [\"524\u00d78=4192\", \"941\u00d77=6587\"],\n  [\"572\u00d73=1716\", \"254\u00d79=2286\"],\n  [\"263\u00d75=1315\", \"296\u00d74=1184\"],\n  [\"846\u00d76=5076\", \"999\u00d77=6993\"],\n  [\"363\u00d78=2904\", \"827\u00d78=6616\"],\n  [\"734\u00d72=1468\", \"102\u00d76=612\"],\n  [\"992\u00d79=8928\", \"242\u00d72=484\"],\n  [\"989\u00d79=8901\", \"866\u00d72=1732\"],\n  [\"747\u00d74=2988\", \"428\u00d72=856\"],\n  [\"432\u00d72=864\", \"883\u00d77=6181\"],\n  [\"895\u00d77=6265\", \"955\u00d78=7640\"],\n  [\"985\u00d76=5910\", \"482\u00d72=964\"],\n  [\"389\u00d77=2723\", \"468\u00d78=3744\"],\n  [\"428\u00d79=3852\", \"315\u00d76=1890\"],\n  [\"386\u00d75=1930\", \"165\u00d72=330\"],\n  [\"751\u00d73=2253\", \"769\u00d75=3845\"],\n  [\"971\u00d79=8739\", \"622\u00d74=2488\"],\n  [\"657\u00d77=4599\", \"524\u00d73=1572\"],\n  [\"856\u00d79=7704\", \"482\u00d75=2410\"],\n  [\"389\u00d73=1167\", \"889\u00d73=2667\"],\n  [\"444\u00d79=3996\", \"634\u00d72=1268\"],\n  [\"261\u00d72=522\", \"825\u00d78=6600\"],\n  [\"475\u00d72=950\", \"525\u00d79=4725\"],\n  [\"873\u00d78=6984\", \"666\u00d79=5994\"],\n  [\"146\u00d75=730\", \"370\u00d79=3330\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"524\u00d78=4192\"; New = \"941\u00d77=6587\" },\n    @{ Old = \"572\u00d73=1716\"; New = \"254\u00d79=2286\" },\n    @{ Old = \"263\u00d75=1315\"; New = \"296\u00d74=1184\" },\n    @{ Old = \"846\u00d76=5076\"; New = \"999\u00d77=6993\" },\n    @{ Old = \"363\u00d78=2904\"; New = \"827\u00d78=6616\" },\n    @{ Old = \"734\u00d72=1468\"; New = \"102\u00d76=612\" },\n    @{ Old = \"992\u00d79=8928\"; New = \"242\u00d72=484\" },\n    @{ Old = \"989\u00d79=8901\"; New = \"866\u00d72=1732\" },\n    @{ Old = \"747\u00d74=2988\"; New = \"428\u00d72=856\" },\n    @{ Old = \"432\u00d72=864\"; New = \"883\u00d77=6181\" },\n    @{ Old = \"895\u00d77=6265\"; New = \"955\u00d78=7640\" },\n    @{ Old = \"985\u00d76=5910\"; New = \"482\u00d72=964\" },\n    @{ Old = \"389\u00d77=2723\"; New = \"468\u00d78=3744\" },\n    @{ Old = \"428\u00d79=3852\"; New = \"315\u00d76=1890\" },\n    @{ Old = \"386\u00d75=1930\"; New = \"165\u00d72=330\" },\n    @{ Old = \"751\u00d73=2253\"; New = \"769\u00d75=3845\" },\n    @{ Old = \"971\u00d79=8739\"; New = \"622\u00d74=2488\" },\n    @{ Old = \"657\u00d77=4599\"; New = \"524\u00d73=1572\" },\n    @{ Old = \"856\u00d79=7704\"; New = \"482\u00d75=2410\" },\n    @{ Old = \"389\u00d73=1167\"; New = \"889\u00d73=2667\" },\n    @{ Old = \"444\u00d79=3996\"; New = \"634\u00d72=1268\" },\n    @{ Old = \"261\u00d72=522\"; New = \"825\u00d78=6600\" },\n    @{ Old = \"475\u00d72=950\"; New = \"525\u00d79=4725\" },\n    @{ Old = \"873\u00d78=6984\"; New = \"666\u00d79=5994\" },\n    @{ Old = \"146\u00d75=730\"; New = \"370\u00d79=3330\" }\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $r.New, $wdReplaceAll) | Out-Null\n}\n"}
